# Fusion des motifs "Se rendre dans une grande surface ou un centre commercial..."
# et "Se rendre dans un centre de proximite..." en un seul motif "Commerce".
#
# motive_id (col C) for both groups becomes 2, and motive_label (col D)
# becomes "Commerce" for rows 2-25 (the two groups that used to carry the
# distinct "grande surface" / "proximite" motive descriptions).
# The now-unused F column amounts that were tied to the old "proximite"
# rows (rows 6-25) are cleared as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4: previously motive_id 2.1 / "grande surface" motive_label.
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 3).Value = 2
    $ws.Cells.Item($r, 4).Value = "Commerce"
}

# Rows 5-25: previously motive_id 2.2 / "proximite" motive_label.
for ($r = 5; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 2
    $ws.Cells.Item($r, 4).Value = "Commerce"
}

# The frequentation values that used to sit on F6:F25 no longer apply once
# these rows are folded into the generic "Commerce" motive.
for ($r = 6; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).ClearContents()
}

# Restore the sheet selection to where the author left off editing.
$ws.Range("B16").Select() | Out-Null
